# Commit: Change "device" to "apparatus" (#13)
# Rename the "Device" worksheet to "Apparatus" and update its descriptive
# text, then make it the active/selected sheet (mirroring the author's
# interactive session state captured in the saved workbook).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Device")
$ws.Name = "Apparatus"
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."

# Make the renamed sheet the active tab/selection, as reflected in the diff
# (activeTab moves to the Apparatus sheet, and its selection becomes A2).
$ws.Activate()
$ws.Range("A2").Select()
